$replacements = @(
    @{ Old = "DTaP [1]"; New = "DTaP " }
    @{ Old = "DTaP-IPV [2]"; New = "DTaP-IPV " }
    @{ Old = "DTaP-Hep B-IPV [4]"; New = "DTaP-Hep B-IPV " }
    @{ Old = "DTaP-IP-HI [4]"; New = "DTaP-IP-HI " }
    @{ Old = "e-IPV [5]"; New = "e-IPV " }
    @{ Old = "Hepatitis A Pediatric [5]"; New = "Hepatitis A Pediatric " }
    @{ Old = "Hepatitis A-Hepatitis B 18 only [3]"; New = "Hepatitis A-Hepatitis B 18 only " }
    @{ Old = "Hepatitis B [5]`nPediatric/Adolescent"; New = "Hepatitis B  Pediatric/Adolescent" }
    @{ Old = "Recombivax`nHB"; New = "Recombivax HB" }
    @{ Old = "Hib [5]"; New = "Hib " }
    @{ Old = "HPV - Human Papillomavirus 9-valent [5]"; New = "HPV - Human Papillomavirus 9-valent " }
    @{ Old = "MENB - Meningococcal Group B [5]"; New = "MENB - Meningococcal Group B " }
    @{ Old = "Meningococcal Conjugate (Groups A, C, Y and W-135) [5]"; New = "Meningococcal Conjugate (Groups A, C, Y and W-135) " }
    @{ Old = "Measles, Mumps and Rubella (MMR) [1]"; New = "Measles, Mumps and Rubella (MMR) " }
    @{ Old = "MMR/Varicella [2]"; New = "MMR/Varicella " }
    @{ Old = "Pneumococcal`n13-valent [5] (Pediatric)"; New = "Pneumococcal 13-valent  (Pediatric)" }
    @{ Old = "Rotavirus, Live, Oral, Pentavalent [5]"; New = "Rotavirus, Live, Oral, Pentavalent " }
    @{ Old = "Rotavirus, Live, Oral, Oral [5]"; New = "Rotavirus, Live, Oral, Oral " }
    @{ Old = "Tetanus and Diphtheria Toxoids [3]"; New = "Tetanus and Diphtheria Toxoids " }
    @{ Old = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis [1]"; New = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis " }
    @{ Old = "Varicella [5]"; New = "Varicella " }
    @{ Old = "Hepatitis A Adult [5]"; New = "Hepatitis A Adult " }
    @{ Old = "Hepatitis A-Hepatitis B Adult [3]"; New = "Hepatitis A-Hepatitis B Adult " }
    @{ Old = "Hepatitis B Adult [5]"; New = "Hepatitis B Adult " }
    @{ Old = "HPV-Human Papillomavirus 9 Valent [5]"; New = "HPV-Human Papillomavirus 9 Valent " }
    @{ Old = "Measles, Mumps,  Rubella [1]"; New = "Measles, Mumps,  Rubella " }
    @{ Old = "Pneumococcal`n13-valent [5]"; New = "Pneumococcal 13-valent " }
    @{ Old = "Influenza [5]`n(Age 6 months and older)"; New = "Influenza  (Age 6 months and older)" }
    @{ Old = "Fluzone`nQuadrivalent"; New = "Fluzone Quadrivalent" }
    @{ Old = "Influenza [5]`n(Age 6-35 months)"; New = "Influenza  (Age 6-35 months)" }
    @{ Old = "Fluzone`nQuadrivalent`nPediatric dose"; New = "Fluzone Quadrivalent Pediatric dose" }
    @{ Old = "Fluarix`nQuadrivalent"; New = "Fluarix Quadrivalent" }
    @{ Old = "FluLaval`nQuadrivalent"; New = "FluLaval Quadrivalent" }
    @{ Old = "Influenza [5]`n(Age 4 years and older)"; New = "Influenza  (Age 4 years and older)" }
    @{ Old = "Influenza [5]`n(Age 6 -35 months)"; New = "Influenza  (Age 6 -35 months)" }
    @{ Old = "Influenza [5]`n(Age 36 months and older)"; New = "Influenza  (Age 36 months and older)" }
    @{ Old = "Influenza [5]`nLive, Intranasal (Age 2-49 years)"; New = "Influenza  Live, Intranasal (Age 2-49 years)" }
    @{ Old = "FluMist`nQuadrivalent"; New = "FluMist Quadrivalent" }
    @{ Old = "Afluria`nQuadrivalent"; New = "Afluria Quadrivalent" }
)

$wb = $excel.ActiveWorkbook

# xlWhole = 1 (match the entire cell content, not a substring - several of the
# old footnote-style labels are literal prefixes of other, unrelated labels,
# e.g. "Fluzone`nQuadrivalent" vs "Fluzone`nQuadrivalent`nPediatric dose", so a
# part/substring match would corrupt unrelated cells).
# xlByRows = 1
$xlWhole = 1
$xlByRows = 1

foreach ($ws in $wb.Worksheets) {
    foreach ($r in $replacements) {
        $null = $ws.Cells.Replace($r.Old, $r.New, $xlWhole, $xlByRows, $false, $false, $false)
    }
}

Write-Host "Done replacing" $replacements.Count "strings across" $wb.Worksheets.Count "sheets"
